$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '62.738.17'
$ws.Range("E2").Value = '  +2.23%  '
$ws.Range("D3").Value = '3.027.84'
$ws.Range("E3").Value = '  +1.91%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.78%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.025.57'
$ws.Range("E8").Value = '  +1.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.519'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("E10").Value = '  +11.08%  '
$ws.Range("E11").Value = '  +5.38%  '
$ws.Range("E12").Value = '  +1.17%  '
$ws.Range("E13").Value = '  +3.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.44%  '
$ws.Range("E15").Value = '  +2.56%  '
$ws.Range("D16").Value = '3.526.14'
$ws.Range("E16").Value = '  +1.82%  '
$ws.Range("D17").Value = '62.640.26'
$ws.Range("E17").Value = '  +2.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").Value = '3.026.63'
$ws.Range("E19").Value = '  +1.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '450.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.692'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.24%  '
$ws.Range("E26").Value = '  +9.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.48%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("E29").Value = '  +3.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.92%  '
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = '  +4.00%  '
$ws.Range("D35").Value = '0.0₃0862'
$ws.Range("E35").Value = '  +11.28%  '
$ws.Range("E36").Value = '  +2.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.88'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.31'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("E42").Value = '  +4.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.297'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +12.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '393.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.68%  '
$ws.Range("E46").Value = '  +0.88%  '
$ws.Range("D47").Value = '2.740.45'
$ws.Range("E47").Value = '  +1.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.62%  '
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.20'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.11%  '
$ws.Range("E51").Value = '  +0.38%  '
